$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "9th Stab - Cosmetic Changes" -------------------------------------
# This weekly watch-sheet keeps one column per reporting date, newest on
# the left (right after the Analyst column A). A new reporting week
# (Jun_17, with an interim Jun_15) is being rolled in, which pushes the
# two existing date columns (old B = Jun_13, old C = Jun_10) two slots to
# the right (-> D, E) to make room.

$lastRow = 27

# Remember each row's current (pre-insert) column-B value so it can be
# replayed into the new column D once column B is free to take the new
# "Jun_17" data.
$oldB = @{}
for ($r = 1; $r -le $lastRow; $r++) {
    $oldB[$r] = $ws.Cells.Item($r, 2).Value2
}

# Insert two fresh columns at C:D. Column B stays put; the old column C
# (Jun_10 data, with all of its cell formatting) slides into column E.
$ws.Range("C:D").EntireColumn.Insert()

# The freshly inserted columns should keep the same display width as the
# column they were cloned from.
$ws.Columns.Item(3).ColumnWidth = 7.14
$ws.Columns.Item(4).ColumnWidth = 7.14
$ws.Columns.Item(5).ColumnWidth = 7.14

# Column D = the old column B data (the Jun_13 column), now shifted over.
for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = $oldB[$r]
}

# Column C = new "Jun_15" column. No ratings/prices changed that week, so
# every analyst just shows "UN" (unchanged) again, same as the header
# row shows the new "Jun_15" date.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
}

# Column B = new "Jun_17" column header + data. Default everyone to "UN"
# (unchanged) first ...
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
}

# ... then the header row dates and the one real rating move this week.
$ws.Cells.Item(1, 2).Value = "Jun_17"
$ws.Cells.Item(1, 3).Value = "Jun_15"

# Row 4 (Morningstar) had a stray "Jun_13" value sitting in its old
# column B instead of "UN" - that same stray value now shows up in both
# of the newly inserted cells (C4 and D4) since it simply rides along as
# the row's column-B content gets replayed/copied across the new layout.
$ws.Cells.Item(4, 3).Value = $oldB[4]

# Row 22 (BidaskClub): a brand-new upgrade was reported on 6/16/2018 -
# put it in the new Jun_17 column and flag it with the same highlight
# fill used elsewhere on this sheet for a freshly-reported rating.
$ws.Cells.Item(22, 2).Value = "6/16/2018,Upgrades,Buy -> Strong-Buy,"
$ws.Cells.Item(22, 2).Interior.ColorIndex = 35
